$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cells($ws, $col, $row1, $row2) {
    $cell1 = $ws.Range($col + $row1)
    $cell2 = $ws.Range($col + $row2)
    $v1 = $cell1.Value()
    $v2 = $cell2.Value()
    $cell1.Value = $v2
    $cell2.Value = $v1
}

function Swap-Row-Columns($ws, $row1, $row2, $cols) {
    foreach ($col in $cols) {
        Swap-Cells $ws $col $row1 $row2
    }
}

# Row 8 / Row 9: taxon details, coordinates and times were swapped between
# the two observation records (rest of the row is identical).
Swap-Row-Columns $ws 8 9 @("A","B","E","F","G","H","Q","R","Z","AB")

# Row 18 / Row 19: only Id, coordinates and times differ between the two
# records, so swap just those.
Swap-Row-Columns $ws 18 19 @("A","Q","R","Z","AB")
